# Generate Report for Handoff
#
# A new localization file (228b4934-faed-44a0-b362-1b99a5cea0b9.md) has been
# handed off. Each worksheet (Overview, zh-cn, de-de) gets a new row for it,
# inserted ahead of the pre-existing 372efa8c-... row, pushing the
# .localization-config row down one more slot.

$wb = $excel.ActiveWorkbook

$newFileId  = "228b4934-faed-44a0-b362-1b99a5cea0b9"
$oldFileId  = "372efa8c-0a36-4dd9-9388-3d2fad088cea"
$newHash    = "c7f65fee7b20e509e2ecb2f5a389c22b40dbfd31"
$oldHash    = "23e37f775862b718012b593c4c81fb24ab6b8f19"
$commitSha  = "92472e064ad4a77888c110c63ba5230aa658d728"

$mdBase       = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha/e2e"
$configUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha/.localization-config"
$zhHandoffSha = "6c621432329d5be56f6a79ad95658ad8386c6bb8"
$deHandoffSha = "ce081b9846c6c20521ca6b6b861794272a948d2a"
$zhBase       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deBase       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$newMdName = "$newFileId.md"
$oldMdName = "$oldFileId.md"
$configName = ".localization-config"

$readyStatus = "Ready for handoff"
$notLocalized = "Not to be localized"
$epoch = "0001-01-01 00:00:00"

# --- Sheet 1: Overview --------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("B2").Value = $readyStatus
$ws1.Range("C2").Value = $readyStatus
$ws1.Range("B3").Value = $readyStatus
$ws1.Range("C3").Value = $readyStatus
$ws1.Range("B4").Value = $notLocalized
$ws1.Range("C4").Value = $notLocalized

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$mdBase/$newMdName", "", "", $newMdName)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$mdBase/$oldMdName", "", "", $oldMdName)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", $configName)

# --- Sheet 2: zh-cn -------------------------------------------------------
$newZhXlfName = "$newFileId.$newHash.zh-cn.xlf"
$oldZhXlfName = "$oldFileId.$oldHash.zh-cn.xlf"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("B2").Value = $readyStatus
$ws2.Range("D2").Value = "2016-02-22 13:58:42"
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = "Include"

$ws2.Range("B3").Value = $readyStatus
$ws2.Range("D3").Value = "2016-02-22 13:58:00"
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = "Include"

$ws2.Range("B4").Value = $notLocalized
$ws2.Range("D4").Value = $epoch
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$mdBase/$newMdName", "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhBase/$newZhXlfName", "", "", $newZhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$mdBase/$oldMdName", "", "", $oldMdName)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhBase/$oldZhXlfName", "", "", $oldZhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", $configName)

# --- Sheet 3: de-de -------------------------------------------------------
$newDeXlfName = "$newFileId.$newHash.de-de.xlf"
$oldDeXlfName = "$oldFileId.$oldHash.de-de.xlf"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("B2").Value = $readyStatus
$ws3.Range("D2").Value = "2016-02-22 13:58:54"
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = "Include"

$ws3.Range("B3").Value = $readyStatus
$ws3.Range("D3").Value = "2016-02-22 13:58:10"
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = "Include"

$ws3.Range("B4").Value = $notLocalized
$ws3.Range("D4").Value = $epoch
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$mdBase/$newMdName", "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deBase/$newDeXlfName", "", "", $newDeXlfName)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$mdBase/$oldMdName", "", "", $oldMdName)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deBase/$oldDeXlfName", "", "", $oldDeXlfName)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", $configName)
